$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("invalidLoginTest")

# Update row 2: password becomes a numeric value instead of "saul1223"
$ws.Range("B2").Value = 123344

# Add new row 4 for "bala" / "bala123" invalid login test data
$ws.Range("A4").Value = "bala"
$ws.Range("B4").Value = "bala123"
$ws.Range("C4").Value = "Invalid credentials"

# Update the selected/active cell to B2 as shown in the saved workbook
$ws.Range("B2").Select()
